# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across the resume body.
#
# Strategy: for each target paragraph, repeatedly use Range.Find.Execute to
# locate each metric substring (search scoped to that paragraph's Range so
# we never bleed into neighboring paragraphs), then set Font.Bold / Font.Color
# on the matched Range. Word's Find automatically splits the run(s) so the
# matched text becomes its own <w:r> with the new formatting while the
# surrounding text keeps the original (unformatted) run(s).

$d = $word.ActiveDocument

# Convert an "RRGGBB" hex string into the BGR-packed integer that
# Word's Font.Color (OLE COLORREF) expects.
function ColorFromHex($hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($bb * 65536) + ($gg * 256) + $rr
}

$metricColor = ColorFromHex('2C3E50')

# Find `text` inside `para`'s own range and apply bold + metric color to it.
function Highlight-Metric($para, $text) {
    $pr = $para.Range
    $searchRange = $d.Range($pr.Start, $pr.End)
    $found = $searchRange.Find.Execute($text, $true, $false, $false, $false, $false, `
                                        $true, 1, $false, '', 0)
    if ($found) {
        $searchRange.Font.Bold = 1
        $searchRange.Font.Color = $metricColor
    }
    return $found
}

# Apply a list of metric substrings (in order) to a single paragraph.
function Highlight-Metrics($para, $metrics) {
    foreach ($m in $metrics) {
        Highlight-Metric $para $m | Out-Null
    }
}

$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text

    if ($t -like '*Discovered systematic race coding errors*') {
        Highlight-Metrics $p @('23%', '64%')
    }
    elseif ($t -like '*Utilized advanced sampling methods*') {
        Highlight-Metrics $p @('±4.2%', '±2.1%', '71%', '87%')
    }
    elseif ($t -like '*Trigonometric algorithm for boundary estimation*') {
        Highlight-Metrics $p @('73.5%', '$4.7M')
    }
    elseif ($t -like '*Built real-time FEC analysis systems*') {
        Highlight-Metrics $p @('$2')
    }
    elseif ($t -like '*Modernized legacy ETL processes*') {
        Highlight-Metrics $p @('57%')
    }
    elseif ($t -like '*Algorithmic innovation: Pioneered trigonometric*') {
        Highlight-Metrics $p @('73.5%')
    }
    elseif ($t -like '*$4.7M savings enabled nonprofit access*') {
        Highlight-Metrics $p @('$4.7M')
    }
    elseif ($t -like '*Platform impact: Built redistricting system serving*') {
        Highlight-Metrics $p @('12,847')
    }
}
